# Apply the commit's changes:
#  1. Swap the table's style GUID on the "Plenary" cash-flow table (slide 16)
#     from the Table_0 custom style to the built-in Medium-Style-2 GUID.
#  2. Swap the two embedded themes: the theme actually driving the deck
#     (Slide Master / Notes Master / Handout Master all resolve to the same
#     live theme object in this host) goes from the "Integral" palette to
#     the stock "Office" palette - i.e. theme2.xml <-> theme1.xml contents.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tableShape = $shape
            break
        }
    }
    if ($tableShape -ne $null) { break }
}

if ($tableShape -ne $null) {
    # Table styles cannot be assigned through a property - ApplyStyle is the
    # supported entry point for swapping a table's styleId.
    $tableShape.Table.ApplyStyle("{C24AE0CB-CBBF-4BD0-8DDB-4D1EBCB9B4A5}", $true)
}

# --- 2. Theme colors --------------------------------------------------------
$theme = $p.SlideMaster.Theme
$cs = $theme.ThemeColorScheme

# Best-effort: keep the visible theme name in sync too (some hosts surface
# this as read-only, so failures here are non-fatal to the rest of the script).
try { $theme.Name = "Office Theme" } catch { }
try { $cs.Name = "Office" } catch { }

# dk1 / lt1 / dk2 / lt2 / accent1-6 / hlink / folHlink, in ThemeColorScheme
# index order, set to the stock "Office" theme palette values.
$cs.Item(1).RGB = 0          # dk1      000000
$cs.Item(2).RGB = 16777215   # lt1      FFFFFF
$cs.Item(3).RGB = 6968388    # dk2      44546A
$cs.Item(4).RGB = 15132391   # lt2      E7E6E6
$cs.Item(5).RGB = 13998939   # accent1  5B9BD5
$cs.Item(6).RGB = 3243501    # accent2  ED7D31
$cs.Item(7).RGB = 10855845   # accent3  A5A5A5
$cs.Item(8).RGB = 49407      # accent4  FFC000
$cs.Item(9).RGB = 12874308   # accent5  4472C4
$cs.Item(10).RGB = 4697456   # accent6  70AD47
$cs.Item(11).RGB = 12673797  # hlink    0563C1
$cs.Item(12).RGB = 7491477   # folHlink 954F72
